$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 232, shifting the existing rows 232-252 down to
# 233-253 (same semantics as right-clicking the row header and choosing
# "Insert").
$ws.Rows.Item(232).Insert()

# Populate the freshly inserted row 232 with the new weekly price record.
$ws.Cells.Item(232, 1).Value = 10
$ws.Cells.Item(232, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(232, 3).Value = "La Araucanía"
$ws.Cells.Item(232, 4).Value = 44578
$ws.Cells.Item(232, 5).Value = 9
$ws.Cells.Item(232, 6).Value = 100114013
$ws.Cells.Item(232, 7).Value = "Zanahoria"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 100
$ws.Cells.Item(232, 11).Value = 8000
$ws.Cells.Item(232, 12).Value = 8000
$ws.Cells.Item(232, 13).Value = 8000
$ws.Cells.Item(232, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(232, 15).Value = "Región del Maule"
$ws.Cells.Item(232, 16).Value = 320
$ws.Cells.Item(232, 17).Value = 25
$ws.Cells.Item(232, 18).Value = "Hortaliza"
